$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 5161.6665
$ws.Range("J7").Value = 7580
$ws.Range("L7").Value = 7580
$ws.Range("N7").Value = -7804
$ws.Range("H9").Value = 174
$ws.Range("I9").Value = 175
$ws.Range("J9").Value = 173
$ws.Range("K9").Value = 175
$ws.Range("L9").Value = 173
$ws.Range("M9").Value = -6
$ws.Range("N9").Value = -511
$ws.Range("H14").Value = 5161.6665
$ws.Range("J14").Value = 7580
$ws.Range("L14").Value = 7580
$ws.Range("N14").Value = -7962
$ws.Range("H29").Value = 4316.5
$ws.Range("I29").Value = 633.3333
$ws.Range("K29").Value = 1899.9999
$ws.Range("M29").Value = -1618.9999
$ws.Range("H38").Value = 2061.4
$ws.Range("I38").Value = 2061.4
$ws.Range("K38").Value = 6184.200000000001
$ws.Range("M38").Value = -5812.200000000001
$ws.Range("H41").Value = 1016.4
$ws.Range("I41").Value = 1434
$ws.Range("J41").Value = 390
$ws.Range("K41").Value = 1434
$ws.Range("L41").Value = 390
$ws.Range("M41").Value = -994
$ws.Range("N41").Value = -1270
$ws.Range("H43").Value = 1111
$ws.Range("I43").Value = 1111
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1111
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = ""
$ws.Range("N43").Value = -1042
$ws.Range("H53").Value = 596.63635
$ws.Range("I53").Value = 456.3
$ws.Range("K53").Value = 456.3
$ws.Range("M53").Value = 180.7
$ws.Range("H58").Value = 1342.3077
$ws.Range("I58").Value = 89
$ws.Range("J58").Value = 2125.625
$ws.Range("K58").Value = 267
$ws.Range("L58").Value = 6376.875
$ws.Range("M58").Value = -117
$ws.Range("N58").Value = -6676.875
$ws.Range("H98").Value = 747.8
$ws.Range("I98").Value = 435
$ws.Range("J98").Value = 1999
$ws.Range("K98").Value = 435
$ws.Range("L98").Value = 1999
$ws.Range("M98").Value = 1063
$ws.Range("N98").Value = -4995
$ws.Range("H116").Value = 11883.571
$ws.Range("I116").Value = 3947.5
$ws.Range("J116").Value = 15058
$ws.Range("K116").Value = 3947.5
$ws.Range("L116").Value = 15058
$ws.Range("M116").Value = -505.5
$ws.Range("N116").Value = -21942
$ws.Range("H122").Value = 747.8
$ws.Range("I122").Value = 435
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 1305
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = 1145
$ws.Range("N122").Value = -10897

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2642.3333
$ws.Range("I122").Value = 1874.75
$ws.Range("K122").Value = 5624.25
$ws.Range("M122").Value = -3174.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 399.75
$ws.Range("J22").Value = 499.5
$ws.Range("L22").Value = 499.5
$ws.Range("N22").Value = -845.5
$ws.Range("H107").Value = 984.53845
$ws.Range("I107").Value = 1084
$ws.Range("J107").Value = 760.75
$ws.Range("K107").Value = 1084
$ws.Range("L107").Value = 760.75
$ws.Range("M107").Value = 836
$ws.Range("N107").Value = -4600.75
$ws.Range("H134").Value = 6505.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3430
$ws.Range("I16").Value = 3430
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3430
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = -3143
$ws.Range("H41").Value = 500
$ws.Range("I41").Value = 500
$ws.Range("K41").Value = 500
$ws.Range("M41").Value = -72
$ws.Range("H60").Value = 10046.5
$ws.Range("I60").Value = 10046.5
$ws.Range("K60").Value = 10046.5
$ws.Range("M60").Value = -9535.5
$ws.Range("H62").Value = 4748.5
$ws.Range("I62").Value = 4748.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4748.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = -4124.5
$ws.Range("H65").Value = 4748.5
$ws.Range("I65").Value = 4748.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 23742.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = -20622.5
$ws.Range("H68").Value = 100000
$ws.Range("J68").Value = 100000
$ws.Range("L68").Value = 100000
$ws.Range("N68").Value = -101498
$ws.Range("H71").Value = 100000
$ws.Range("J71").Value = 100000
$ws.Range("L71").Value = 300000
$ws.Range("N71").Value = -307488
$ws.Range("H99").Value = 2778.8
$ws.Range("I99").Value = 1736
$ws.Range("K99").Value = 1736
$ws.Range("M99").Value = -238
$ws.Range("H113").Value = 3430
$ws.Range("I113").Value = 3430
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3430
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -1260
$ws.Range("H126").Value = 2778.8
$ws.Range("I126").Value = 1736
$ws.Range("K126").Value = 5208
$ws.Range("M126").Value = -2738

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 4000
$ws.Range("J104").Value = 4000
$ws.Range("L104").Value = 12000
$ws.Range("N104").Value = -17242
$ws.Range("H126").Value = 5400
$ws.Range("I126").Value = 5400
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 16200
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = -11260
$ws.Range("H140").Value = 3998.5
$ws.Range("I140").Value = 3998.5
$ws.Range("K140").Value = 11995.5
$ws.Range("M140").Value = -6815.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6432.25
$ws.Range("I70").Value = 6139.6
$ws.Range("J70").Value = 6641.2856
$ws.Range("K70").Value = 6139.6
$ws.Range("L70").Value = 6641.2856
$ws.Range("M70").Value = -5869.6
$ws.Range("N70").Value = -7181.2856
$ws.Range("H73").Value = 6432.25
$ws.Range("I73").Value = 6139.6
$ws.Range("J73").Value = 6641.2856
$ws.Range("K73").Value = 6139.6
$ws.Range("L73").Value = 6641.2856
$ws.Range("M73").Value = -5203.6
$ws.Range("N73").Value = -8513.285599999999
$ws.Range("H122").Value = 2082.6667
$ws.Range("I122").Value = 2082.6667
$ws.Range("K122").Value = 6248.000100000001
$ws.Range("M122").Value = -3798.000100000001
$ws.Range("H126").Value = 3633.25
$ws.Range("I126").Value = 4011
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 12033
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -9563
$ws.Range("N126").Value = -12440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 90001
$ws.Range("I18").Value = 90001
$ws.Range("K18").Value = 90001
$ws.Range("M18").Value = -89829
$ws.Range("H40").Value = 2199.5
$ws.Range("I40").Value = 2199.5
$ws.Range("K40").Value = 2199.5
$ws.Range("M40").Value = -2063.5
$ws.Range("H46").Value = 1142.7142
$ws.Range("I46").Value = 1499.6666
$ws.Range("J46").Value = 875
$ws.Range("K46").Value = 1499.6666
$ws.Range("L46").Value = 875
$ws.Range("M46").Value = -1311.6666
$ws.Range("N46").Value = -1251
$ws.Range("H68").Value = 30500
$ws.Range("I68").Value = 30500
$ws.Range("K68").Value = 30500
$ws.Range("M68").Value = -29751
$ws.Range("H71").Value = 30500
$ws.Range("I71").Value = 30500
$ws.Range("K71").Value = 152500
$ws.Range("M71").Value = -148756
$ws.Range("H95").Value = 39999
$ws.Range("J95").Value = 39999
$ws.Range("L95").Value = 39999
$ws.Range("N95").Value = -45491
$ws.Range("H122").Value = 4797.25
$ws.Range("I122").Value = 2729.6667
$ws.Range("K122").Value = 8189.000100000001
$ws.Range("M122").Value = -5739.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 50000
$ws.Range("J75").Value = 50000
$ws.Range("L75").Value = 50000
$ws.Range("N75").Value = -51872
$ws.Range("H78").Value = 50000
$ws.Range("J78").Value = 50000
$ws.Range("L78").Value = 150000
$ws.Range("N78").Value = -159360
$ws.Range("H107").Value = 686.1111
$ws.Range("I107").Value = 686.1111
$ws.Range("K107").Value = 2058.3333
$ws.Range("M107").Value = -138.3332999999998
$ws.Range("H122").Value = 12566.444
$ws.Range("I122").Value = 1849.75
$ws.Range("J122").Value = 21139.8
$ws.Range("K122").Value = 5549.25
$ws.Range("L122").Value = 63419.39999999999
$ws.Range("M122").Value = -3099.25
$ws.Range("N122").Value = -68319.39999999999
